$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column E (old E/F shift right to G/H, etc.)
$ws.Range("E1:F1").EntireColumn.Insert()

# New header cells
$ws.Range("E1").Value = "CO Amount"
$ws.Range("F1").Value = "Cash-in-Hand"

# New row-2 values (the "type" row)
$ws.Range("E2").Value = "Numeric for C/O Refi & Debt Consolidation. Not used for others"
$ws.Range("F2").Value = "Numeric for C/O Refi & Debt Consolidation. Not used for others"

# Match the column width of column D
$ws.Range("E1:F1").EntireColumn.ColumnWidth = $ws.Range("D1").EntireColumn.ColumnWidth

# Update the active selection
$ws.Range("F2").Select()
